# Update the "Förändrad" (Changed) date column (column C) for rows 2 through 34
# from 2023-09-15 (serial 45184) to 2023-09-16 (serial 45185).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 34; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45184) {
        $cell.Value2 = 45185
    }
}

